# Add a new row (row 80) of data to each of the 4 worksheets, mirroring
# the existing last row (row 79) but with an updated timestamp.

$wb = $excel.ActiveWorkbook

# Data for the new row 80 on each worksheet, in sheet order.
$rowsData = @(
    @{ A = 45866.43819444445; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x3C"; E = "0x14"; F = 380; G = [double]"7.598631275147109e+23"; H = 316; I = 14 },
    @{ A = 45866.43819444445; B = "0x01,0x7c"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x3C"; E = "0xe";  F = 380; G = [double]"5.68432987514711e+23";  H = 316; I = 14 },
    @{ A = 45866.43819444445; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x78"; E = "0x7";  F = 130; G = [double]"5.68631262647114e+23";  H = 120; I = 7 },
    @{ A = 45866.43819444445; B = "0x00,0x82"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x77"; E = "0x3";  F = 130; G = [double]"9.85046333984776e+23";  H = 119; I = 3 }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $data = $rowsData[$i]
    $newRow = 80

    # Column A keeps the same date number format as the row above it.
    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
